# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, shifting the existing "Late" / "heading" / "Outstanding" columns
# one place to the right, and leave that sheet active/selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make this the active sheet (drives workbook.xml's activeTab and this
# sheet's tabSelected attribute).
$ws.Activate()

# Insert a new blank column at N; everything from N onward shifts right.
$ws.Columns("N").Insert()

# The freshly inserted column picks up the width of its neighbour (column M).
$ws.Columns("N").ColumnWidth = 10.2

# Reflect the new selection left behind on this sheet.
$ws.Range("J18").Select() | Out-Null
